$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking" row): Right 5 -> 4, Wrong -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total" row): Right 40 -> 32, Wrong -3 -> -6, Max text "40 / 140" -> "26 / 112"
$ws.Range("B12").Value = 32
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "26 / 112"
